# Engine_5/vets_test.xlsx edit: "Completed remove_whitespace script in jupyter"
#
# The author's whitespace-cleanup script touched a handful of Last/First
# Name cells, replacing an odd (non-breaking-space) leading character with
# plain spaces -- but, per the commit message, it over-corrected on a few
# rows and left double leading spaces instead of trimming them. It also
# moved the viewport/selection before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content fixes -----------------------------------------------
# (order matches the new shared-string append order of the reference edit)

# Row 59: BUCKO, " THEODORE W" (leading NBSP) -> "  THEODORE W" (two spaces)
$ws.Range("B59").Value = "  THEODORE W"

# Row 79: "MARTEL" -> " MARTEL" (gains a leading space)
$ws.Range("A79").Value = " MARTEL"

# Row 74: WUJICK, " EDWARD W" -> "  EDWARD W" (extra leading space)
$ws.Range("B74").Value = "  EDWARD W"

# Row 65: "CAMARA" -> " CAMARA" (gains a leading space)
$ws.Range("A65").Value = " CAMARA"

# --- View state ----------------------------------------------------------
# Move the active selection to A65 (was E73), scrolling the window up a bit.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$ws.Range("A65").Select()
